$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.821.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.352.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.75%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.38%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.29"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.30%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.16%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -8.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.708.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.321.14"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.823.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "77.09"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.68"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.41%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.00%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.40"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.78"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.42"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.98%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.94%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.75%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0891"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.96%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.54%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0360"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.70%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.35"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.231"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "114.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.86"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.30%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.17"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.52"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.33%  "
